$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 42649.656493055554
$ws.Range("A8").NumberFormat = "m/d/yy h:mm"

$ws.Range("B8").Value = $true

$ws.Range("C8").Value = 10043.799999999999
$ws.Range("D8").Value = 9993.33
$ws.Range("E8").Value = 18.829999999999998
$ws.Range("F8").Value = 19.02

$ws.Range("G8").Value = $false
$ws.Range("G8").NumberFormat = "m/d/yy h:mm"

$ws.Range("H8").Value = 1.01

$ws.Range("I8").Value = $false
